# The commit removes the "section_summary" column (and its now-unused
# shared-string values) from the sheet. In the original workbook this was
# column D (section_title, section_dialogue, section_time_stamp,
# section_summary, episode_title, episode_date). Deleting the entire D
# column shifts episode_title/episode_date from E/F into D/E, matching
# the diff (dimension A1:F8 -> A1:E8, cols 5/6 -> 4/5, etc.) and drops the
# section_summary header + its 7 paragraph strings from sharedStrings.xml
# since they're no longer referenced anywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Select()
$ws.Columns("D").Delete()
